# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The 'K' column (spreadsheet column G) values are recalculated and rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(4,2,1,3,2,0,2,0,3,3,1,2,2,1,1,0,1,0,1,2,1,1,1,0,1,0,1,1,1,0,2,1,0,1,0,0,1,2,2,3,3,2,1,1,0,0,1,0,4,2,0,4,2,2,1,0,1,2,2,1,0,1,2,0,0,1)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
